# Auto-generated PowerShell/Word-COM script implementing the
# "Sprawozdanie - score board" commit:
#   1) splits the "Po srodku ramki..." sentence, inserting the word " gry"
#   2) adds a new paragraph about the game-time display (right side of the
#      score frame), in List-Paragraph style
#   3) adds a blank spacer paragraph (List-Paragraph style)
#   4) adds a new "5.3 Przyciski" heading paragraph (List-Paragraph style,
#      bigger font)
#   5) fills what used to be the final, empty paragraph of the document with
#      the new "Nad tablica wynikow... przyciskow." paragraph (tab, tab,
#      text) -- its own paragraph formatting is left untouched.

$d = $word.ActiveDocument

# --- Step 1: locate the *whole original sentence* of paragraph A via Find.
#     The match's Start/End give us the exact text range to rebuild (this
#     excludes the paragraph's two leading <w:tab/> runs and its trailing
#     paragraph mark). We then grow the range left by 1 character so it also
#     swallows the *second* tab (which lives in the very same <w:r> as the
#     text, right before <w:t>), keeping the first tab-only run untouched. ---
$searchHit = $d.Content
$ok = $searchHit.Find.Execute("Po środku ramki jest wyświetlany aktualny wynik, który jest uaktualniany przy każdym ruchu, który zmienił stan planszy. ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok) {
    throw "Could not locate the 'Po srodku ramki...' sentence"
}
$foundStart = $searchHit.Start
$foundEnd = $searchHit.End

$replaceA = $d.Range($foundStart - 1, $foundEnd)
$xmlA = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii="Calibri Light" w:hAnsi="Calibri Light" w:cs="Calibri Light"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:tab/><w:t>Po środku ramki jest wyświetlany aktualny wynik</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri Light" w:hAnsi="Calibri Light" w:cs="Calibri Light"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:t xml:space="preserve"> gry</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri Light" w:hAnsi="Calibri Light" w:cs="Calibri Light"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:t xml:space="preserve">, który jest uaktualniany przy każdym ruchu, który zmienił stan planszy. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$replaceA.InsertXML($xmlA)

# --- Step 2: insert the three brand-new paragraphs just before the very
#     last paragraph of the document body (found via Content.End, which is
#     robust against any phantom/duplicate trailing paragraph entries) ---
$endPos = $d.Content.End
$lastParaRange = $d.Range($endPos - 1, $endPos).Paragraphs.Item(1).Range
$insertPoint = $d.Range($lastParaRange.Start, $lastParaRange.Start)
$xmlNewParas = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Akapitzlist"/><w:ind w:left="624"/><w:contextualSpacing w:val="0"/><w:rPr><w:rFonts w:ascii="Calibri Light" w:hAnsi="Calibri Light" w:cs="Calibri Light"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri Light" w:hAnsi="Calibri Light" w:cs="Calibri Light"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri Light" w:hAnsi="Calibri Light" w:cs="Calibri Light"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:tab/><w:t>Po prawej stronie jest wyświetlany czas gry</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri Light" w:hAnsi="Calibri Light" w:cs="Calibri Light"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:t xml:space="preserve"> w formacie HH:MM:SS. Czas jest uaktualniany co 100 milisekund. </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Akapitzlist"/><w:ind w:left="624"/><w:contextualSpacing w:val="0"/><w:rPr><w:rFonts w:ascii="Calibri Light" w:hAnsi="Calibri Light" w:cs="Calibri Light"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="Akapitzlist"/><w:ind w:left="624"/><w:contextualSpacing w:val="0"/><w:rPr><w:rFonts w:ascii="Calibri Light" w:hAnsi="Calibri Light" w:cs="Calibri Light"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri Light" w:hAnsi="Calibri Light" w:cs="Calibri Light"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>5.3 Przyciski</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$insertPoint.InsertXML($xmlNewParas)

# --- Step 3: fill the (still empty) final paragraph with the new
#     "Nad tablica wynikow..." buttons text; its own pPr is untouched ---
$endPos2 = $d.Content.End
$finalParaRange = $d.Range($endPos2 - 1, $endPos2).Paragraphs.Item(1).Range
$fillPoint = $d.Range($finalParaRange.Start, $finalParaRange.Start)
$xmlE = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii="Calibri Light" w:hAnsi="Calibri Light" w:cs="Calibri Light"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri Light" w:hAnsi="Calibri Light" w:cs="Calibri Light"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:tab/><w:t xml:space="preserve">Nad tablicą wyników znajduje się pięć przycisków. Na środku ekranu znajduje się przycisk z napisem ‘ZRESTATRUJ GRĘ’  naciśnięcie go powoduje zrestartowanie gry wraz ze zrestartowaniem wyniku i czasu gry. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$fillPoint.InsertXML($xmlE)

Write-Output "Edit applied."
